$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-6 (columns B and C) and add new rows 7-15 (columns A, B, C)
# with the newly imported translation strings.

$ws.Range("B2").Value = 'Gestión de Corpus Documentales'
$ws.Range("C2").Value = 'Documental Corpus Management'

$ws.Range("B3").Value = 'Idioma'
$ws.Range("C3").Value = 'Language'

$ws.Range("B4").Value = 'Inicio'
$ws.Range("C4").Value = 'Home'

$ws.Range("B5").Value = 'Añadir documentos'
$ws.Range("C5").Value = 'Add documents'

$ws.Range("B6").Value = 'Visualizar datos'
$ws.Range("C6").Value = 'Visualize data'

$ws.Range("A7").Value = 'menu.intro.title'
$ws.Range("B7").Value = 'Bienvenido al servicio web de visualización gráfica de conjuntos de datos'
$ws.Range("C7").Value = 'Welcome to your dataset graphic visualization web service.'

$ws.Range("A8").Value = 'menu.intro.first_paragraph'
$ws.Range("B8").Value = 'El mundo de la información evoluciona. Nos vemos obligados cada día a adaptar nuestra tecnología a la oleada de cantidades de datos cada vez más extensas. A la vez que aumenta la necesidad de generar sistemas que aseguren su protección y capacidad de salvaguarda, es necesario el desarrollo de nuevos servicios que nos permitan conocer de una manera auxiliar y sencilla las características de los conjuntos de información que poseemos.'
$ws.Range("C8").Value = 'The world of global information is evolving. Everyday, we''re forced to adapt our technological environment for upcoming waves of extensive data amounts. At the same time it''s necessary to generate system to protect information, it''s necessary to develop new services to access our documents'' most important characteristics, in the simplest of ways.'

$ws.Range("A9").Value = 'menu.intro.second_paragraph'
$ws.Range("B9").Value = 'Por ello, este proyecto, desarrollado inicialmente como un Trabajo de Fin de Grado, pretende poner de manifiesto las utilidades llevadas por tecnologías web que nos permiten analizar nuestros propios corpus documentales, por medio de gráficas interactivas.'
$ws.Range("C9").Value = 'Because of that, this project - developed initially as a final university degree project - tries to manifest the ultimate web techonologies features which permit us to analyze our own documental corpuses using interactive graphs.'

$ws.Range("A10").Value = 'menu.indexList.title'
$ws.Range("B10").Value = 'Consulta la lista actual de conjuntos en línea'
$ws.Range("C10").Value = 'Check the current online index list'

$ws.Range("A11").Value = 'menu.indexList.subtitle'
$ws.Range("B11").Value = 'Selecciona un corpus para empezar'
$ws.Range("C11").Value = 'Select a corpus to begin'

$ws.Range("A12").Value = 'menu.indexList.inputPlaceholder'
$ws.Range("B12").Value = '… o crea un nuevo índice'
$ws.Range("C12").Value = '… or create a new one'

$ws.Range("A13").Value = 'menu.indexList.create'
$ws.Range("B13").Value = 'Crear'
$ws.Range("C13").Value = 'Create'

$ws.Range("A14").Value = 'menu.indexList.documents'
$ws.Range("B14").Value = 'documentos'
$ws.Range("C14").Value = 'documents'

$ws.Range("A15").Value = 'menu.indexList.creating'
$ws.Range("B15").Value = 'Creando el índice…'
$ws.Range("C15").Value = 'Creating index…'

# Update the active cell selection to A12 (matches the post-edit workbook state)
$ws.Range("A12").Select()
